$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '69.653.94'
$ws.Range("E2").Value = '  -0.08%  '
$ws.Range("D3").Value = '3.676.29'
$ws.Range("E3").Value = '  -0.66%  '
$ws.Range("E4").Value = '  +0.03%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '650.70'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -3.75%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '160.99'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -0.85%  '
$ws.Range("E7").Value = '  +0.05%  '
$ws.Range("E8").Value = '  -0.07%  '
$ws.Range("E9").Value = '  -2.13%  '
$ws.Range("E10").Value = '  +0.64%  '
$ws.Range("E11").Value = '  -0.96%  '
$ws.Range("E12").Value = '  -2.47%  '
$ws.Range("D13").Value = '4.296.92'
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '32.67'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -0.77%  '
$ws.Range("D15").Value = '3.669.47'
$ws.Range("E15").Value = '  -0.47%  '
$ws.Range("D16").Value = '69.714.98'
$ws.Range("E16").Value = '  +0.03%  '
$ws.Range("E17").Value = '  +0.86%  '
$ws.Range("E18").Value = '  -0.05%  '
$ws.Range("E19").Value = '  -1.63%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '10.32'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +4.95%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '471.14'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -0.46%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.653'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -0.36%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '79.87'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -0.79%  '
$ws.Range("D24").Value = '3.823.85'
$ws.Range("E24").Value = '  -0.65%  '
$ws.Range("E25").Value = '  +0.02%  '
$ws.Range("E26").Value = '  -1.54%  '
$ws.Range("E27").Value = '  +0.85%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '8.77'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -4.09%  '
$ws.Range("E29").Value = '  -2.36%  '
$ws.Range("E30").Value = '  -3.32%  '
$ws.Range("B31").Value = 'Binance-PegBSC-USD'
$ws.Range("C31").Value = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.00'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -0.06%  '
$ws.Range("B32").Value = 'ImmutableX'
$ws.Range("C32").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '1.99'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -1.91%  '
$ws.Range("E33").Value = '  +0.81%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '26.69'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -0.85%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '6.41'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -3.11%  '
$ws.Range("D36").Value = '3.671.34'
$ws.Range("E36").Value = '  -0.47%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '8.36'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -1.74%  '
$ws.Range("B39").Value = 'Monero'
$ws.Range("C39").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '178.72'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +4.97%  '
$ws.Range("B40").Value = 'Filecoin'
$ws.Range("C40").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '5.86'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -5.28%  '
$ws.Range("E41").Value = '  -0.04%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.0893'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -1.62%  '
$ws.Range("E43").Value = '  -2.37%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.928'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -1.83%  '
$ws.Range("B45").Value = 'InjectiveProtocol'
$ws.Range("C45").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '29.15'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +4.18%  '
$ws.Range("B46").Value = 'OKB'
$ws.Range("C46").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '46.67'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -0.69%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '2.77'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -1.14%  '
$ws.Range("B48").Value = 'FLOKI'
$ws.Range("C48").Value = 'https://coinranking.com/coin/fmHk13Rqw+floki-floki'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.000268'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -4.56%  '
$ws.Range("B49").Value = 'Cosmos'
$ws.Range("C49").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '7.86'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -0.79%  '
$ws.Range("B50").Value = 'SuiNetwork'
$ws.Range("C50").Value = 'https://coinranking.com/coin/3xJluUMvp+suinetwork-sui'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '1.05'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -5.70%  '
$ws.Range("B51").Value = 'ONDO'
$ws.Range("C51").Value = 'https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '1.23'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -5.38%  '
